# Updated cryptos list data (Price and Volume(1h) columns), matching the
# upstream GitHub Actions scraper commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.409.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.10%  '

$ws.Range("D3").Value = "'1.946.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.53%  '

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.39%  '

$ws.Range("D5").Value = "'325.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").Value = "'0.4626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.47%  '

$ws.Range("E8").Value = '  -0.52%  '

$ws.Range("D9").Value = "'46.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.70%  '

$ws.Range("D10").Value = "'0.07829"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.39%  '

$ws.Range("D11").Value = "'0.9793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.14%  '

$ws.Range("D12").Value = "'22.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.65%  '

$ws.Range("D13").Value = "'1.946.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.54%  '

$ws.Range("D14").Value = "'7.076"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("D15").Value = "'5.749"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.39%  '

$ws.Range("D16").Value = "'0.07051"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("D17").Value = "'86.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").Value = "'0.000009815"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.17%  '

$ws.Range("D20").Value = "'17.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").Value = "'29.433.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("D23").Value = "'5.464"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.74%  '

$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("D25").Value = "'2.167.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.16%  '

$ws.Range("D26").Value = "'2.097"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = "'157.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.84%  '

$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("D29").Value = "'5.757"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.02%  '

$ws.Range("D30").Value = "'118.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").Value = "'1.862"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.72%  '

$ws.Range("D32").Value = "'0.09372"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.62%  '

$ws.Range("E33").Value = '  -3.80%  '

$ws.Range("E34").Value = '  -1.04%  '

$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("D36").Value = "'3.127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.76%  '

$ws.Range("D37").Value = "'0.05765"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.58%  '

$ws.Range("D38").Value = "'1.158"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.18%  '

$ws.Range("D39").Value = "'0.02084"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.20%  '

$ws.Range("D40").Value = "'7.693"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.17%  '

$ws.Range("D41").Value = "'0.5657"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.30%  '

$ws.Range("D42").Value = "'0.1783"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("D43").Value = "'9.446"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.87%  '

$ws.Range("D44").Value = "'0.000002860"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +36.23%  '

$ws.Range("D45").Value = "'2.733"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.12%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.5288"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.89%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'11.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.04%  '

$ws.Range("D48").Value = "'2.099"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.04%  '

$ws.Range("D49").Value = "'0.06871"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.92%  '

$ws.Range("D50").Value = "'1.814"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.86%  '

$ws.Range("D51").Value = "'111.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.18%  '
